# rerun corona results with larger ds
#
# This script rewrites the "negative" (A-column) and "positive" (J-column)
# anchor-word tables on Sheet1 to reflect a rerun against a larger dataset.
# The negative-word table shrinks from 8 words to 5 words (rows 3-7), and
# the positive-word table grows from 29 words to 31 words (rows 3-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# New data for the "negative" table (columns A-H), rows 3-7.
# Columns: row, name, anchor score, type occurences, total occurences,
#          +%, -%, both (bool), normal
# ----------------------------------------------------------------------
$aRows = @(
    @(3, "crude", 0.7941176470588235, 27, 27, 0, 1, $false, 7),
    @(4, "crisis", 0.5856164383561644, 171, 171, 0, 1, $false, 121),
    @(5, "fraud", 0.5833333333333334, 21, 21, 0, 1, $false, 15),
    @(6, "sc", 0.2063492063492063, 39, 39, 0, 1, $false, 150),
    @(7, "panic", 0.1705426356589147, 88, 88, 0, 1, $false, 428)
)

# ----------------------------------------------------------------------
# New data for the "positive" table (columns J-Q), rows 3-33.
# Columns: row, name, anchor score, type occurences, total occurences,
#          +%, -%, both (bool), normal
# ----------------------------------------------------------------------
$jRows = @(
    @(3, "happy", 1, 26, 26, 1, 0, $false, 0),
    @(4, "love", 0.9347826086956522, 43, 43, 1, 0, $false, 3),
    @(5, "best", 0.9152542372881356, 54, 54, 1, 0, $false, 5),
    @(6, "interesting", 0.8787878787878788, 29, 29, 1, 0, $false, 4),
    @(7, "great", 0.8660714285714286, 97, 97, 1, 0, $false, 15),
    @(8, "thank", 0.8203125, 105, 105, 1, 0, $false, 23),
    @(9, "nice", 0.8148148148148148, 22, 22, 1, 0, $false, 5),
    @(10, "positive", 0.7931034482758621, 46, 46, 1, 0, $false, 12),
    @(11, "thanks", 0.7926829268292683, 65, 65, 1, 0, $false, 17),
    @(12, "free", 0.7583333333333333, 91, 91, 1, 0, $false, 29),
    @(13, "special", 0.7222222222222222, 26, 26, 1, 0, $false, 10),
    @(14, "safety", 0.7058823529411765, 36, 36, 1, 0, $false, 15),
    @(15, "safe", 0.6971830985915493, 99, 99, 1, 0, $false, 43),
    @(16, "confidence", 0.6944444444444444, 25, 25, 1, 0, $false, 11),
    @(17, "support", 0.6792452830188679, 72, 72, 1, 0, $false, 34),
    @(18, "good", 0.6625, 106, 106, 1, 0, $false, 54),
    @(19, "fresh", 0.6041666666666666, 29, 29, 1, 0, $false, 19),
    @(20, "heroes", 0.574468085106383, 27, 27, 1, 0, $false, 20),
    @(21, "relief", 0.5600000000000001, 28, 28, 1, 0, $false, 22),
    @(22, "well", 0.5425531914893617, 51, 51, 1, 0, $false, 43),
    @(23, "better", 0.5238095238095238, 33, 33, 1, 0, $false, 30),
    @(24, "hand", 0.4882506527415144, 187, 187, 1, 0, $false, 196),
    @(25, "like", 0.4558823529411765, 155, 155, 1, 0, $false, 185),
    @(26, "care", 0.4269662921348314, 38, 38, 1, 0, $false, 51),
    @(27, "help", 0.423728813559322, 125, 125, 1, 0, $false, 170),
    @(28, "sure", 0.390625, 25, 25, 1, 0, $false, 39),
    @(29, "hope", 0.3692307692307693, 24, 24, 1, 0, $false, 41),
    @(30, "protect", 0.3561643835616438, 26, 26, 1, 0, $false, 47),
    @(31, "increase", 0.3076923076923077, 24, 24, 1, 0, $false, 54),
    @(32, "please", 0.301255230125523, 72, 72, 1, 0, $false, 167),
    @(33, "19", 0.00979934671021932, 21, 21, 1, 0, $false, 2122)
)

# ----------------------------------------------------------------------
# 1. Refresh the A1 title cell (still "negative", kept explicit so the
#    shared-string table is rebuilt the way a fresh export would do it).
# ----------------------------------------------------------------------
$ws.Range("A1").Value = "negative"

# ----------------------------------------------------------------------
# 2. Write the new negative-word rows (3-7), reusing the existing
#    row-3 formatting for the "name" column (A) so style stays intact.
# ----------------------------------------------------------------------
foreach ($row in $aRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

# ----------------------------------------------------------------------
# 3. Remove the now-unused negative-word rows (8-10) entirely so no
#    stray cells remain (the table shrank from 8 to 5 data rows).
# ----------------------------------------------------------------------
$ws.Range("A8:H10").Clear()

# ----------------------------------------------------------------------
# 4. Write the new/updated positive-word rows (3-31 already exist,
#    32-33 are brand-new). For the brand-new rows, first clone the
#    formatting of an existing "name" cell (column J) so the bold
#    header-style border/alignment carries over.
# ----------------------------------------------------------------------
foreach ($row in $jRows) {
    $r = $row[0]

    if ($r -gt 31) {
        $ws.Range("J7").Copy()
        $ws.Range("J" + $r).PasteSpecial(-4122)
    }

    if ($row[1] -eq "19") {
        # Force this shared-string to be stored as text, not a number,
        # since the source word happens to look numeric.
        $ws.Cells.Item($r, 10).NumberFormat = "@"
        $ws.Cells.Item($r, 10).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 10).Value = $row[1]
    }

    $ws.Cells.Item($r, 11).Value = $row[2]
    $ws.Cells.Item($r, 12).Value = $row[3]
    $ws.Cells.Item($r, 13).Value = $row[4]
    $ws.Cells.Item($r, 14).Value = $row[5]
    $ws.Cells.Item($r, 15).Value = $row[6]
    $ws.Cells.Item($r, 16).Value = $row[7]
    $ws.Cells.Item($r, 17).Value = $row[8]
}
